$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.684.85"
$ws.Range("E2").Value = "  -3.44%  "
$ws.Range("D3").Value = "3.350.31"
$ws.Range("E3").Value = "  -1.06%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'569.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.68%  "
$ws.Range("D6").Value = "'133.64"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.97%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.350.41"
$ws.Range("E8").Value = "  -1.11%  "
$ws.Range("E9").Value = "  -0.41%  "
$ws.Range("D10").Value = "'7.56"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.45%  "
$ws.Range("E11").Value = "  +2.02%  "
$ws.Range("E12").Value = "  +2.57%  "
$ws.Range("D13").Value = "3.929.06"
$ws.Range("E13").Value = "  -0.55%  "
$ws.Range("E14").Value = "  +1.39%  "
$ws.Range("D15").Value = "3.357.26"
$ws.Range("E15").Value = "  -0.62%  "
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("E17").Value = "  +1.82%  "
$ws.Range("D18").Value = "60.825.04"
$ws.Range("E18").Value = "  -3.20%  "
$ws.Range("D19").Value = "'13.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.49%  "
$ws.Range("D20").Value = "'9.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.47%  "
$ws.Range("E21").Value = "  +1.87%  "
$ws.Range("D22").Value = "'371.82"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.62%  "
$ws.Range("D23").Value = "'0.570"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.66%  "
$ws.Range("D24").Value = "3.487.17"
$ws.Range("E24").Value = "  -0.83%  "
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").Value = "'70.35"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.70%  "
$ws.Range("E27").Value = "  +9.58%  "
$ws.Range("D28").Value = "'1.65"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +17.51%  "
$ws.Range("D29").Value = "'7.61"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +8.80%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("D31").Value = "'8.05"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.70%  "
$ws.Range("E32").Value = "  +0.65%  "
$ws.Range("D33").Value = "'0.153"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.71%  "
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").Value = "3.382.09"
$ws.Range("E35").Value = "  -0.86%  "
$ws.Range("D36").Value = "'23.27"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.12%  "
$ws.Range("D37").Value = "'5.51"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.74%  "
$ws.Range("D38").Value = "'6.86"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.87%  "
$ws.Range("E39").Value = "  +2.80%  "
$ws.Range("D40").Value = "'162.28"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.28%  "
$ws.Range("E41").Value = "  +3.15%  "
$ws.Range("E42").Value = "  +0.20%  "
$ws.Range("D43").Value = "'41.19"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.44%  "
$ws.Range("B44").Value = "ONDO"
$ws.Range("C44").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D44").Value = "'1.20"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.61%  "
$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").Value = "'4.37"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.09%  "
$ws.Range("D46").Value = "'0.754"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.66%  "
$ws.Range("D47").Value = "'1.58"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.14%  "
$ws.Range("D48").Value = "'6.92"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.75%  "
$ws.Range("D49").Value = "'22.71"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.91%  "
$ws.Range("D50").Value = "'23.20"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +13.46%  "
$ws.Range("E51").Value = "  +13.05%  "
